$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new work-log entry on row 35
$ws.Range("A35").Value = 44061
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = "Backendin testaamista ja confirmaatio sähköpostin lähettämistä"

# Row height grows to fit the wrapped text, same as the other multi-line rows
$ws.Rows.Item(35).RowHeight = 30

# Move the active selection to where the user clicked next
$ws.Range("H35").Select()
